$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount value in T2
$ws.Range("T2").Value = 163769

# Move the active selection to R12 (also shifts the viewport towards it)
[void]$ws.Range("R12").Select()
